# Adição de novos arquivos
# Adds two new "ADM" rows (Antonio ADM / Pedro ADM) to the roster sheet,
# normalizes a handful of leftover "filled" cell styles in column A,
# and sets the sheet up for printing (A4 / portrait).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) A handful of student-name cells (rows 6, 8, 16, 18, 27) still carried
#    a leftover "apply fill" flavour of the border style from earlier
#    edits even though no fill color was ever set. Normalize them to the
#    same plain bordered style used by every other name cell.
# ---------------------------------------------------------------------
$normalizeRows = 6, 8, 16, 18, 27
foreach ($r in $normalizeRows) {
    $ws.Range("A$r").Interior.Pattern = -4142   # xlPatternNone
}

# ---------------------------------------------------------------------
# 2) Append the two new rows of data at the bottom of the roster.
# ---------------------------------------------------------------------
$ws.Range("A37").Value = "Antonio ADM"
$ws.Range("B37").NumberFormat = "00000"
$ws.Range("B37").HorizontalAlignment = -4108    # xlCenter
$ws.Range("B37").Value = 554384356465
$ws.Range("C37").NumberFormat = "00000"
$ws.Range("C37").HorizontalAlignment = -4108    # xlCenter
$ws.Range("C37").Value = 554384356465

$ws.Range("A38").Value = "Pedro ADM"
$ws.Range("B38").NumberFormat = "00000"
$ws.Range("B38").HorizontalAlignment = -4108    # xlCenter
$ws.Range("B38").Value = 5543996440402
$ws.Range("C38").NumberFormat = "00000"
$ws.Range("C38").HorizontalAlignment = -4108    # xlCenter
$ws.Range("C38").Value = 5543996440402

# ---------------------------------------------------------------------
# 3) Reflect the new rows in the window: scroll down a bit and select
#    the two freshly-entered rows (as if the user had just finished
#    typing them in).
# ---------------------------------------------------------------------
$ws.Range("A34").Select()
$ws.Rows("37:38").Select()

# ---------------------------------------------------------------------
# 4) Page setup for printing (A4, portrait) - picked up the one time the
#    sheet was sent to Print Preview.
# ---------------------------------------------------------------------
$ps = $ws.PageSetup
$ps.PaperSize = 9     # xlPaperA4
$ps.Orientation = 1   # xlPortrait
